# Insert a new weekly price record as row 70 (pushing existing rows 70-141
# down to 71-142) on the "Poroto granado" Hortaliza price sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 70..141 down to 71..142, leaving a blank row 70 to populate.
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new weekly record.
$ws.Cells.Item(70, 1).Value = 8
$ws.Cells.Item(70, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(70, 3).Value = "Coquimbo"
$ws.Cells.Item(70, 4).Value = 45280
$ws.Cells.Item(70, 5).Value = 4
$ws.Cells.Item(70, 6).Value = 100112030
$ws.Cells.Item(70, 7).Value = "Poroto granado"
$ws.Cells.Item(70, 8).Value = "Sin especificar"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 400
$ws.Cells.Item(70, 11).Value = 39000
$ws.Cells.Item(70, 12).Value = 40000
$ws.Cells.Item(70, 13).Value = 39500
$ws.Cells.Item(70, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(70, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(70, 16).Value = 1580
$ws.Cells.Item(70, 17).Value = 25
$ws.Cells.Item(70, 18).Value = "Hortaliza"
